# Append a new row (31) to Sheet1 with the next day's gold price data,
# reusing the formatting (borders / wrap text) from the previous row (30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (30) down to the new row (31)
$ws.Range("A30:B30").Copy()
$ws.Range("A31:B31").PasteSpecial(-4122)

# Set the new values
$ws.Range("A31").Value = "22-10-2025"
$ws.Range("B31").Value = "The price of gold in India today is ₹12,720 per gram for 24 karat gold, ₹11,660 per gram for 22 karat gold and ₹9,540 per gram for 18 karat gold (also called 999 gold)."
